# Fix OCR-garbled species names in the port-level landings table (Table35).
# Each correction below replaces a mis-OCR'd species label in column B with
# its corrected reading, matching the same row's port (column A) context.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value  = "Yellnwfin tuna"
$ws.Range("B5").Value  = "Albacore  "
$ws.Range("B7").Value  = "Rockfish"
$ws.Range("B15").Value = "Rex sole"
$ws.Range("B16").Value = "Lingcod  "
$ws.Range("B23").Value = "English sole  "
$ws.Range("B25").Value = "Rockfish"
$ws.Range("B36").Value = "Rockfish"
$ws.Range("B44").Value = "Albacore  "
$ws.Range("B67").Value = "Giant Pacific oyster"
$ws.Range("B66").Value = "Eastern oyster"

# Restore the cursor/selection position left by the author after the edits.
$ws.Range("B67").Select()
